$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Adam10"
$ws.Range("C2").Value = "Epha3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 51.37659299999999
$ws.Range("H2").Value = 154.129779
$ws.Range("I2").Value = 0.6482346823708168
$ws.Range("J2").Value = 0.6482346823708167
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.06698166666666668
$ws.Range("N2").Value = 0.200945
$ws.Range("O2").Value = 0.003012576978541733
$ws.Range("P2").Value = 0.003012576978541732
$ws.Range("Q2").Value = 3.441289826795
$ws.Range("R2").Value = 30.971608441155
$ws.Range("S2").Value = 0.001952856880802635
$ws.Range("T2").Value = 0.001952856880802635

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Adam10"
$ws.Range("C3").Value = "Epha3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 51.37659299999999
$ws.Range("H3").Value = 154.129779
$ws.Range("I3").Value = 0.6482346823708168
$ws.Range("J3").Value = 0.6482346823708167
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 22.03620333333333
$ws.Range("N3").Value = 66.10861
$ws.Range("O3").Value = 0.9911034191912899
$ws.Range("P3").Value = 0.9911034191912899
$ws.Range("Q3").Value = 1132.14504992191
$ws.Range("R3").Value = 10189.30544929719
$ws.Range("S3").Value = 0.6424676101360963
$ws.Range("T3").Value = 0.6424676101360962

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Adam10"
$ws.Range("C4").Value = "Epha3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 51.37659299999999
$ws.Range("H4").Value = 154.129779
$ws.Range("I4").Value = 0.6482346823708168
$ws.Range("J4").Value = 0.6482346823708167
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.130825
$ws.Range("N4").Value = 0.392475
$ws.Range("O4").Value = 0.005884003830168287
$ws.Range("P4").Value = 0.005884003830168287
$ws.Range("Q4").Value = 6.721342779224999
$ws.Range("R4").Value = 60.492085013025
$ws.Range("S4").Value = 0.003814215353917809
$ws.Range("T4").Value = 0.003814215353917808

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Adam10"
$ws.Range("C5").Value = "Epha3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 16.943638
$ws.Range("H5").Value = 50.830914
$ws.Range("I5").Value = 0.2137832260916193
$ws.Range("J5").Value = 0.2137832260916192
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.06698166666666668
$ws.Range("N5").Value = 0.200945
$ws.Range("O5").Value = 0.003012576978541733
$ws.Range("P5").Value = 0.003012576978541732
$ws.Range("Q5").Value = 1.134913112636667
$ws.Range("R5").Value = 10.21421801373
$ws.Range("S5").Value = 0.0006440384253219945
$ws.Range("T5").Value = 0.0006440384253219943

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Adam10"
$ws.Range("C6").Value = "Epha3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 16.943638
$ws.Range("H6").Value = 50.830914
$ws.Range("I6").Value = 0.2137832260916193
$ws.Range("J6").Value = 0.2137832260916192
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 22.03620333333333
$ws.Range("N6").Value = 66.10861
$ws.Range("O6").Value = 0.9911034191912899
$ws.Range("P6").Value = 0.9911034191912899
$ws.Range("Q6").Value = 373.3734521743933
$ws.Range("R6").Value = 3360.36106956954
$ws.Range("S6").Value = 0.2118812863451484
$ws.Range("T6").Value = 0.2118812863451484

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Adam10"
$ws.Range("C7").Value = "Epha3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 16.943638
$ws.Range("H7").Value = 50.830914
$ws.Range("I7").Value = 0.2137832260916193
$ws.Range("J7").Value = 0.2137832260916192
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.130825
$ws.Range("N7").Value = 0.392475
$ws.Range("O7").Value = 0.005884003830168287
$ws.Range("P7").Value = 0.005884003830168287
$ws.Range("Q7").Value = 2.21665144135
$ws.Range("R7").Value = 19.94986297215
$ws.Range("S7").Value = 0.00125790132114882
$ws.Range("T7").Value = 0.00125790132114882

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Adam10"
$ws.Range("C8").Value = "Epha3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 10.935931
$ws.Range("H8").Value = 32.807793
$ws.Range("I8").Value = 0.137982091537564
$ws.Range("J8").Value = 0.137982091537564
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.06698166666666668
$ws.Range("N8").Value = 0.200945
$ws.Range("O8").Value = 0.003012576978541733
$ws.Range("P8").Value = 0.003012576978541732
$ws.Range("Q8").Value = 0.7325068849316666
$ws.Range("R8").Value = 6.592561964384999
$ws.Range("S8").Value = 0.0004156816724171033
$ws.Range("T8").Value = 0.0004156816724171032

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Adam10"
$ws.Range("C9").Value = "Epha3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 10.935931
$ws.Range("H9").Value = 32.807793
$ws.Range("I9").Value = 0.137982091537564
$ws.Range("J9").Value = 0.137982091537564
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 22.03620333333333
$ws.Range("N9").Value = 66.10861
$ws.Range("O9").Value = 0.9911034191912899
$ws.Range("P9").Value = 0.9911034191912899
$ws.Range("Q9").Value = 240.9863991553033
$ws.Range("R9").Value = 2168.87759239773
$ws.Range("S9").Value = 0.1367545227100452
$ws.Range("T9").Value = 0.1367545227100452

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Adam10"
$ws.Range("C10").Value = "Epha3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 10.935931
$ws.Range("H10").Value = 32.807793
$ws.Range("I10").Value = 0.137982091537564
$ws.Range("J10").Value = 0.137982091537564
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.130825
$ws.Range("N10").Value = 0.392475
$ws.Range("O10").Value = 0.005884003830168287
$ws.Range("P10").Value = 0.005884003830168287
$ws.Range("Q10").Value = 1.430693173075
$ws.Range("R10").Value = 12.876238557675
$ws.Range("S10").Value = 0.0008118871551016576
$ws.Range("T10").Value = 0.0008118871551016576
